$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('#cleof', 'Cleof'),
    @('#gaze', 'Gaze'),
    @('#daraie', 'Daraie'),
    @('#darayd', 'DARAYD'),
    @('#larde', 'Larde'),
    @('#neel', 'Neel'),
    @('#sidon', 'SIDON'),
    @('#bruserb', 'Bruserb'),
    @('#ian-soet', 'Ian Soet'),
    @('#darayd', 'Darayd'),
    @('#larde', 'LARDE'),
    @('#dian', 'Dian'),
    @('#sidon', 'Sidon'),
    @('#agasi', 'Agasi'),
    @('#briant', 'Briant'),
    @('#daraid', 'Daraid'),
    @('#floris', 'Floris'),
    @('#diaen', 'Diaen'),
    @('#diane', 'Diane'),
    @('#lard', 'Lard'),
    @('#garay', 'Garay'),
    @('#flori', 'FLORI'),
    @('#sidon,', 'SIDON,'),
    @('#arlang', 'Arlang'),
    @('#larden', 'LARDEN'),
    @('#gund', 'Gund'),
    @('#keyser', 'Keyser'),
    @('#flori', 'Flori'),
    @('#sidonia', 'Sidonia'),
    @('#ages', 'Ages'),
    @('#flor', 'Flor'),
    @('#garai', 'Garai'),
    @('#diana', 'Diana'),
    @('#griet-haecx', 'Griet haecx'),
    @('#bruser', 'Bruser'),
    @('#arlan', 'Arlan'),
    @('#garai', 'GARAI'),
    @('#lard', 'LARD'),
    @('#ian-soete', 'Ian Soete'),
    @('#sidoni', 'Sidoni'),
    @('#dian', 'DIAN'),
    @('#darai', 'Darai'),
    @('#agesi', 'Agesi'),
    @('#galth', 'Galth'),
    @('#flora', 'Flora'),
    @('#darai', 'DARAI'),
    @('#tryn', 'Tryn'),
    @('#garaie', 'Garaie'),
    @('#galta', 'Galta'),
    @('#daray', 'DARAY'),
    @('#cle', 'CLE'),
    @('#garaye', 'Garaye'),
    @('#ian', 'Ian'),
    @('#larden', 'Larden'),
    @('#cleof', 'CLEOF'),
    @('#daray', 'Daray'),
    @('#dar', 'DAR'),
    @('#bruse', 'Bruse'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

# Column D (is_prefered) no longer carries any "x" marker - clear rows 2-27
# which previously held one (rows 28-59 were already blank).
$ws.Range("D2:D27").Value = ""
